$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C2").Value = "Arica y Parinacota"
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D2").Value = 44454
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = 100112043
$ws.Range("G2").Value = "Pepino dulce"
$ws.Range("H2").Value = "Cultivar IV Región"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 160
$ws.Range("K2").Value = 19000
$ws.Range("L2").Value = 20000
$ws.Range("M2").Value = 19500
$ws.Range("N2").Value = "$/bandeja 18 kilos"
$ws.Range("O2").Value = "Provincia de Limarí"
$ws.Range("P2").Value = 1083
$ws.Range("Q2").Value = 18
$ws.Range("R2").Value = "Hortaliza"

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C3").Value = "Arica y Parinacota"
$ws.Range("D3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D3").Value = 44398
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = 100112043
$ws.Range("G3").Value = "Pepino dulce"
$ws.Range("H3").Value = "Cultivar IV Región"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 17000
$ws.Range("L3").Value = 18000
$ws.Range("M3").Value = 17500
$ws.Range("N3").Value = "$/bandeja 18 kilos"
$ws.Range("O3").Value = "Provincia de Limarí"
$ws.Range("P3").Value = 972
$ws.Range("Q3").Value = 18
$ws.Range("R3").Value = "Hortaliza"

# Row 4
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D4").Value = 44398
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = 100112043
$ws.Range("G4").Value = "Pepino dulce"
$ws.Range("H4").Value = "Cultivar IV Región"
$ws.Range("I4").Value = "Segunda"
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 15000
$ws.Range("L4").Value = 16000
$ws.Range("M4").Value = 15500
$ws.Range("N4").Value = "$/bandeja 18 kilos"
$ws.Range("O4").Value = "Provincia de Limarí"
$ws.Range("P4").Value = 861
$ws.Range("Q4").Value = 18
$ws.Range("R4").Value = "Hortaliza"

# Row 5
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C5").Value = "Arica y Parinacota"
$ws.Range("D5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D5").Value = 44412
$ws.Range("E5").Value = 15
$ws.Range("F5").Value = 100112043
$ws.Range("G5").Value = "Pepino dulce"
$ws.Range("H5").Value = "Cultivar IV Región"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 17000
$ws.Range("L5").Value = 18000
$ws.Range("M5").Value = 17500
$ws.Range("N5").Value = "$/bandeja 18 kilos"
$ws.Range("O5").Value = "Provincia de Limarí"
$ws.Range("P5").Value = 972
$ws.Range("Q5").Value = 18
$ws.Range("R5").Value = "Hortaliza"

# Row 6
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C6").Value = "Arica y Parinacota"
$ws.Range("D6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D6").Value = 44377
$ws.Range("E6").Value = 15
$ws.Range("F6").Value = 100112043
$ws.Range("G6").Value = "Pepino dulce"
$ws.Range("H6").Value = "Cultivar IV Región"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 17000
$ws.Range("L6").Value = 18000
$ws.Range("M6").Value = 17600
$ws.Range("N6").Value = "$/bandeja 18 kilos"
$ws.Range("O6").Value = "Provincia de Limarí"
$ws.Range("P6").Value = 978
$ws.Range("Q6").Value = 18
$ws.Range("R6").Value = "Hortaliza"

# Row 7
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C7").Value = "Arica y Parinacota"
$ws.Range("D7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D7").Value = 44435
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = 100112043
$ws.Range("G7").Value = "Pepino dulce"
$ws.Range("H7").Value = "Cultivar IV Región"
$ws.Range("I7").Value = "Segunda"
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 17000
$ws.Range("L7").Value = 18000
$ws.Range("M7").Value = 17500
$ws.Range("N7").Value = "$/bandeja 18 kilos"
$ws.Range("O7").Value = "Provincia de Limarí"
$ws.Range("P7").Value = 972
$ws.Range("Q7").Value = 18
$ws.Range("R7").Value = "Hortaliza"

# Row 8
$ws.Range("A8").Value = 1
$ws.Range("B8").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C8").Value = "Arica y Parinacota"
$ws.Range("D8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D8").Value = 44435
$ws.Range("E8").Value = 15
$ws.Range("F8").Value = 100112043
$ws.Range("G8").Value = "Pepino dulce"
$ws.Range("H8").Value = "Cultivar IV Región"
$ws.Range("I8").Value = "Tercera"
$ws.Range("J8").Value = 120
$ws.Range("K8").Value = 14000
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = 14500
$ws.Range("N8").Value = "$/bandeja 18 kilos"
$ws.Range("O8").Value = "Provincia de Limarí"
$ws.Range("P8").Value = 806
$ws.Range("Q8").Value = 18
$ws.Range("R8").Value = "Hortaliza"

# Row 9
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C9").Value = "Arica y Parinacota"
$ws.Range("D9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D9").Value = 44391
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = 100112043
$ws.Range("G9").Value = "Pepino dulce"
$ws.Range("H9").Value = "Cultivar IV Región"
$ws.Range("I9").Value = "Segunda"
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 15000
$ws.Range("L9").Value = 16000
$ws.Range("M9").Value = 15500
$ws.Range("N9").Value = "$/bandeja 18 kilos"
$ws.Range("O9").Value = "Provincia de Limarí"
$ws.Range("P9").Value = 861
$ws.Range("Q9").Value = 18
$ws.Range("R9").Value = "Hortaliza"

# Row 10
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C10").Value = "Arica y Parinacota"
$ws.Range("D10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D10").Value = 44405
$ws.Range("E10").Value = 15
$ws.Range("F10").Value = 100112043
$ws.Range("G10").Value = "Pepino dulce"
$ws.Range("H10").Value = "Cultivar IV Región"
$ws.Range("I10").Value = "Segunda"
$ws.Range("J10").Value = 140
$ws.Range("K10").Value = 17000
$ws.Range("L10").Value = 18000
$ws.Range("M10").Value = 17500
$ws.Range("N10").Value = "$/bandeja 18 kilos"
$ws.Range("O10").Value = "Provincia de Limarí"
$ws.Range("P10").Value = 972
$ws.Range("Q10").Value = 18
$ws.Range("R10").Value = "Hortaliza"

# Row 11
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C11").Value = "Arica y Parinacota"
$ws.Range("D11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D11").Value = 44433
$ws.Range("E11").Value = 15
$ws.Range("F11").Value = 100112043
$ws.Range("G11").Value = "Pepino dulce"
$ws.Range("H11").Value = "Cultivar IV Región"
$ws.Range("I11").Value = "Segunda"
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 17000
$ws.Range("L11").Value = 18000
$ws.Range("M11").Value = 17500
$ws.Range("N11").Value = "$/bandeja 18 kilos"
$ws.Range("O11").Value = "Provincia de Limarí"
$ws.Range("P11").Value = 972
$ws.Range("Q11").Value = 18
$ws.Range("R11").Value = "Hortaliza"

# Row 12
$ws.Range("A12").Value = 1
$ws.Range("B12").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C12").Value = "Arica y Parinacota"
$ws.Range("D12").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D12").Value = 44433
$ws.Range("E12").Value = 15
$ws.Range("F12").Value = 100112043
$ws.Range("G12").Value = "Pepino dulce"
$ws.Range("H12").Value = "Cultivar IV Región"
$ws.Range("I12").Value = "Tercera"
$ws.Range("J12").Value = 120
$ws.Range("K12").Value = 14000
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = 14500
$ws.Range("N12").Value = "$/bandeja 18 kilos"
$ws.Range("O12").Value = "Provincia de Limarí"
$ws.Range("P12").Value = 806
$ws.Range("Q12").Value = 18
$ws.Range("R12").Value = "Hortaliza"

# Row 13
$ws.Range("A13").Value = 1
$ws.Range("B13").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C13").Value = "Arica y Parinacota"
$ws.Range("D13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D13").Value = 44221
$ws.Range("E13").Value = 15
$ws.Range("F13").Value = 100112043
$ws.Range("G13").Value = "Pepino dulce"
$ws.Range("H13").Value = "Cultivar XV región"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 140
$ws.Range("K13").Value = 5000
$ws.Range("L13").Value = 6000
$ws.Range("M13").Value = 5500
$ws.Range("N13").Value = "$/caja 10 kilos"
$ws.Range("O13").Value = "Región de Arica y Parinacota"
$ws.Range("P13").Value = 550
$ws.Range("Q13").Value = 10
$ws.Range("R13").Value = "Hortaliza"

# Row 14
$ws.Range("A14").Value = 1
$ws.Range("B14").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C14").Value = "Arica y Parinacota"
$ws.Range("D14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D14").Value = 44363
$ws.Range("E14").Value = 15
$ws.Range("F14").Value = 100112043
$ws.Range("G14").Value = "Pepino dulce"
$ws.Range("H14").Value = "Cultivar IV Región"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 140
$ws.Range("K14").Value = 14000
$ws.Range("L14").Value = 15000
$ws.Range("M14").Value = 14500
$ws.Range("N14").Value = "$/bandeja 18 kilos"
$ws.Range("O14").Value = "Provincia de Limarí"
$ws.Range("P14").Value = 806
$ws.Range("Q14").Value = 18
$ws.Range("R14").Value = "Hortaliza"

# Row 15
$ws.Range("A15").Value = 1
$ws.Range("B15").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C15").Value = "Arica y Parinacota"
$ws.Range("D15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D15").Value = 44211
$ws.Range("E15").Value = 15
$ws.Range("F15").Value = 100112043
$ws.Range("G15").Value = "Pepino dulce"
$ws.Range("H15").Value = "Cultivar XV región"
$ws.Range("I15").Value = "Segunda"
$ws.Range("J15").Value = 140
$ws.Range("K15").Value = 4500
$ws.Range("L15").Value = 5000
$ws.Range("M15").Value = 4750
$ws.Range("N15").Value = "$/caja 10 kilos"
$ws.Range("O15").Value = "Región de Arica y Parinacota"
$ws.Range("P15").Value = 475
$ws.Range("Q15").Value = 10
$ws.Range("R15").Value = "Hortaliza"
